$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4958.3335
$ws.Cells.Item(64, 10).Value = 4944.4443
$ws.Cells.Item(64, 12).Value = 4944.4443
$ws.Cells.Item(64, 14).Value = -5440.4443

$ws.Cells.Item(67, 8).Value = 4958.3335
$ws.Cells.Item(67, 10).Value = 4944.4443
$ws.Cells.Item(67, 12).Value = 4944.4443
$ws.Cells.Item(67, 14).Value = -6660.4443

$ws.Cells.Item(70, 8).Value = 4433
$ws.Cells.Item(70, 9).Value = 2902
$ws.Cells.Item(70, 10).Value = 4624.375
$ws.Cells.Item(70, 11).Value = 8706
$ws.Cells.Item(70, 12).Value = 13873.125
$ws.Cells.Item(70, 13).Value = -8436
$ws.Cells.Item(70, 14).Value = -14413.125

$ws.Cells.Item(73, 8).Value = 4433
$ws.Cells.Item(73, 9).Value = 2902
$ws.Cells.Item(73, 10).Value = 4624.375
$ws.Cells.Item(73, 11).Value = 8706
$ws.Cells.Item(73, 12).Value = 13873.125
$ws.Cells.Item(73, 13).Value = -7770
$ws.Cells.Item(73, 14).Value = -15745.125

$ws.Cells.Item(88, 8).Value = 13174.5
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 13174.5
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 13).Value = 13174.5
$ws.Cells.Item(88, 14).Value = -13986.5
$ws.Cells.Item(88, 12).ClearContents()

$ws.Cells.Item(91, 8).Value = 13174.5
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 13174.5
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 13).Value = 13174.5
$ws.Cells.Item(91, 14).Value = -15982.5
$ws.Cells.Item(91, 12).ClearContents()

$ws.Cells.Item(98, 8).Value = 2048.4
$ws.Cells.Item(98, 9).Value = 1212.5714
$ws.Cells.Item(98, 10).Value = 3998.6667
$ws.Cells.Item(98, 11).Value = 1212.5714
$ws.Cells.Item(98, 12).Value = 3998.6667
$ws.Cells.Item(98, 13).Value = 285.4286
$ws.Cells.Item(98, 14).Value = -6994.6667

$ws.Cells.Item(100, 8).Value = 3071.4285
$ws.Cells.Item(100, 9).Value = 2960
$ws.Cells.Item(100, 10).Value = 3350
$ws.Cells.Item(100, 11).Value = 2960
$ws.Cells.Item(100, 12).Value = 3350
$ws.Cells.Item(100, 13).Value = -2419
$ws.Cells.Item(100, 14).Value = -4432

$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 14).Value = 0
$ws.Cells.Item(110, 12).ClearContents()

$ws.Cells.Item(122, 8).Value = 2048.4
$ws.Cells.Item(122, 9).Value = 1212.5714
$ws.Cells.Item(122, 10).Value = 3998.6667
$ws.Cells.Item(122, 11).Value = 3637.7142
$ws.Cells.Item(122, 12).Value = 11996.0001
$ws.Cells.Item(122, 13).Value = -1187.7142
$ws.Cells.Item(122, 14).Value = -16896.0001

$ws.Cells.Item(141, 8).Value = 3438.0833
$ws.Cells.Item(141, 9).Value = 3438.0833
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 10314.2499
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = -5134.249899999999
$ws.Cells.Item(141, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2394.6667
$ws.Cells.Item(32, 9).Value = 2209.018
$ws.Cells.Item(32, 11).Value = 2209.018
$ws.Cells.Item(32, 13).Value = -1922.018

$ws.Cells.Item(45, 8).Value = 2115.4546
$ws.Cells.Item(45, 9).Value = 1862.3334
$ws.Cells.Item(45, 11).Value = 1862.3334
$ws.Cells.Item(45, 13).Value = -1485.3334

$ws.Cells.Item(74, 8).Value = 1075.5333
$ws.Cells.Item(74, 9).Value = 1075.5333
$ws.Cells.Item(74, 11).Value = 1075.5333
$ws.Cells.Item(74, 13).Value = -201.5333000000001

$ws.Cells.Item(77, 8).Value = 1075.5333
$ws.Cells.Item(77, 9).Value = 1075.5333
$ws.Cells.Item(77, 11).Value = 5377.6665
$ws.Cells.Item(77, 13).Value = -1009.6665

$ws.Cells.Item(97, 8).Value = 904
$ws.Cells.Item(97, 9).Value = 822.3333
$ws.Cells.Item(97, 11).Value = 822.3333
$ws.Cells.Item(97, 13).Value = -326.3333

$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 14).Value = 0
$ws.Cells.Item(109, 12).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(95, 8).Value = 26541.334
$ws.Cells.Item(95, 10).Value = 26541.334
$ws.Cells.Item(95, 12).Value = 26541.334
$ws.Cells.Item(95, 14).Value = -32033.334

$ws.Cells.Item(134, 8).Value = 10433.728
$ws.Cells.Item(134, 9).Value = 11077.1
$ws.Cells.Item(134, 11).Value = 33231.3
$ws.Cells.Item(134, 13).Value = -30696.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 333379.94
$ws.Cells.Item(7, 9).Value = 41.583332
$ws.Cells.Item(7, 11).Value = 41.583332
$ws.Cells.Item(7, 13).Value = 71.416668

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2976.8572
$ws.Cells.Item(34, 9).Value = 1200
$ws.Cells.Item(34, 10).Value = 3113.5386
$ws.Cells.Item(34, 11).Value = 3600
$ws.Cells.Item(34, 12).Value = 9340.6158
$ws.Cells.Item(34, 13).Value = -3516
$ws.Cells.Item(34, 14).Value = -9508.6158

$ws.Cells.Item(39, 8).Value = 8094.5
$ws.Cells.Item(39, 9).Value = 2500
$ws.Cells.Item(39, 10).Value = 9959.333000000001
$ws.Cells.Item(39, 11).Value = 7500
$ws.Cells.Item(39, 12).Value = 29877.999
$ws.Cells.Item(39, 13).Value = -7206
$ws.Cells.Item(39, 14).Value = -30465.999

$ws.Cells.Item(55, 8).Value = 1897.8422
$ws.Cells.Item(55, 10).Value = 2190.9333
$ws.Cells.Item(55, 12).Value = 6572.7999
$ws.Cells.Item(55, 14).Value = -6926.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 3059
$ws.Cells.Item(9, 9).Value = 732.8
$ws.Cells.Item(9, 10).Value = 6936
$ws.Cells.Item(9, 11).Value = 732.8
$ws.Cells.Item(9, 12).Value = 6936
$ws.Cells.Item(9, 13).Value = -562.8
$ws.Cells.Item(9, 14).Value = -7276

$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 14).Value = 0
$ws.Cells.Item(64, 12).ClearContents()

$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 14).Value = 0
$ws.Cells.Item(67, 12).ClearContents()

$ws.Cells.Item(97, 8).Value = 794.4
$ws.Cells.Item(97, 9).Value = 772.1111
$ws.Cells.Item(97, 10).Value = 995
$ws.Cells.Item(97, 11).Value = 772.1111
$ws.Cells.Item(97, 12).Value = 995
$ws.Cells.Item(97, 13).Value = -276.1111
$ws.Cells.Item(97, 14).Value = -1987

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 44254.223
$ws.Cells.Item(7, 9).Value = 43916.234
$ws.Cells.Item(7, 11).Value = 43916.234
$ws.Cells.Item(7, 13).Value = -43804.234

$ws.Cells.Item(126, 8).Value = 44254.223
$ws.Cells.Item(126, 9).Value = 43916.234
$ws.Cells.Item(126, 11).Value = 131748.702
$ws.Cells.Item(126, 13).Value = -129278.702

$ws.Cells.Item(132, 8).Value = 3101.1177
$ws.Cells.Item(132, 9).Value = 2476.5833
$ws.Cells.Item(132, 11).Value = 7429.749899999999
$ws.Cells.Item(132, 13).Value = -4899.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 13124.5
$ws.Cells.Item(74, 9).Value = 13832
$ws.Cells.Item(74, 10).Value = 12770.75
$ws.Cells.Item(74, 11).Value = 13832
$ws.Cells.Item(74, 12).Value = 12770.75
$ws.Cells.Item(74, 13).Value = -12896
$ws.Cells.Item(74, 14).Value = -14642.75

$ws.Cells.Item(77, 8).Value = 13124.5
$ws.Cells.Item(77, 9).Value = 13832
$ws.Cells.Item(77, 10).Value = 12770.75
$ws.Cells.Item(77, 11).Value = 41496
$ws.Cells.Item(77, 12).Value = 38312.25
$ws.Cells.Item(77, 13).Value = -36816
$ws.Cells.Item(77, 14).Value = -47672.25
